$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.201037526130676
$ws.Range("B1").Value = 2.062352180480957
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.030741214752197
$ws.Range("E1").Value = 1.206578135490417
